$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 value (typo'd test data string changed from "5#$#(&" to "35#$#(&")
$ws.Range("D2").Value = "35#$#(&"

# G2 used to hold a numeric phone-like value (111222) formatted as an integer.
# It now holds a text value ("1219234") using the same text style as the rest
# of the row (numFmtId 49 / "@" -> applied via NumberFormat "@").
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1219234"

# Move the active selection from G2 to A2
$ws.Range("A2").Select()
